# Omaha_Cal_Info_GP05MOAS-PG514_00001.xlsx
# "NUTNR Cal Coef for GPGs" - Added NUTNR Cal Coefficients for GPG 514, 528, 566

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Insert 6 new rows for the additional NUTNR calibration coefficients ---
# Row 17 already holds the GP05MOAS-PG514-05-NUTNRM000 / "No calibration
# coefficient" placeholder row; six new rows (18-23) are inserted right
# after it so that the old rows 18-22 move down to 24-28.
$ws.Rows("18:23").Insert()

# --- Fill the newly inserted rows with the repeated Ref Des / barcode /
#     serial / deployment / sensor barcode / sensor serial info, mirroring
#     row 17, then set the per-row calibration coefficient name + value.

$refDes       = "GP05MOAS-PG514-05-NUTNRM000"
$moorBarcode  = "A01222"
$moorSerial   = 514
$deployNum    = 1
$sensBarcode  = "N00652"
$sensSerial   = 357

for ($r = 17; $r -le 23; $r++) {
    $ws.Range("A$r").Value2 = $refDes
    $ws.Range("B$r").Value2 = $moorBarcode
    $ws.Range("C$r").Value2 = $moorSerial
    $ws.Range("D$r").Value2 = $deployNum
    $ws.Range("E$r").Value2 = $sensBarcode
    $ws.Range("F$r").Value2 = $sensSerial
}

# Calibration coefficient *names* (column G) are entered first, for every
# new row in order, so their shared-string ids land contiguously ...
$ws.Range("G17").Value2 = "CC_cal_temp"
$ws.Range("G18").Value2 = "CC_wl"
$ws.Range("G19").Value2 = "CC_eno3"
$ws.Range("G20").Value2 = "CC_eswa"
$ws.Range("G21").Value2 = "CC_di"
$ws.Range("G22").Value2 = "CC_lower_wavelength_limit_for_spectra_fit"
$ws.Range("G23").Value2 = "CC_upper_wavelength_limit_for_spectra_fit"

# ... and only then the coefficient *values* (column H), so the long array
# literals get their own shared-string ids after the names.
$ws.Range("H17").Value2 = 19.98
$ws.Range("H18").Value2 = "[189.61, 190.41, 191.21, 192, 192.8, 193.6, 194.4, 195.2, 196, 196.8, 197.6, 198.4, 199.2, 200, 200.81, 201.61, 202.41, 203.21, 204.02, 204.82, 205.62, 206.42, 207.23, 208.03, 208.84, 209.64, 210.44, 211.25, 212.05, 212.86, 213.67, 214.47, 215.28, 216.08, 216.89, 217.7, 218.5, 219.31, 220.12, 220.92, 221.73, 222.54, 223.35, 224.16, 224.96, 225.77, 226.58, 227.39, 228.2, 229.01, 229.82, 230.63, 231.44, 232.25, 233.06, 233.87, 234.68, 235.49, 236.3, 237.11, 237.93, 238.74, 239.55, 240.36, 241.17, 241.99, 242.8, 243.61, 244.42, 245.24, 246.05, 246.86, 247.68, 248.49, 249.3, 250.12, 250.93, 251.75, 252.56, 253.37, 254.19, 255, 255.82, 256.63, 257.45, 258.26, 259.08, 259.89, 260.71, 261.52, 262.34, 263.16, 263.97, 264.79, 265.6, 266.42, 267.24, 268.05, 268.87, 269.69, 270.5, 271.32, 272.14, 272.95, 273.77, 274.59, 275.41, 276.22, 277.04, 277.86, 278.68, 279.49, 280.31, 281.13, 281.95, 282.77, 283.58, 284.4, 285.22, 286.04, 286.86, 287.67, 288.49, 289.31, 290.13, 290.95, 291.77, 292.59, 293.4, 294.22, 295.04, 295.86, 296.68, 297.5, 298.32, 299.14, 299.96, 300.77, 301.59, 302.41, 303.23, 304.05, 304.87, 305.69, 306.51, 307.33, 308.15, 308.96, 309.78, 310.6, 311.42, 312.24, 313.06, 313.88, 314.7, 315.52, 316.34, 317.16, 317.97, 318.79, 319.61, 320.43, 321.25, 322.07, 322.89, 323.71, 324.53, 325.35, 326.16, 326.98, 327.8, 328.62, 329.44, 330.26, 331.08, 331.89, 332.71, 333.53, 334.35, 335.17, 335.99, 336.8, 337.62, 338.44, 339.26, 340.08, 340.89, 341.71, 342.53, 343.35, 344.16, 344.98, 345.8, 346.62, 347.43, 348.25, 349.07, 349.89, 350.7, 351.52, 352.34, 353.15, 353.97, 354.79, 355.6, 356.42, 357.23, 358.05, 358.87, 359.68, 360.5, 361.31, 362.13, 362.94, 363.76, 364.57, 365.39, 366.2, 367.02, 367.83, 368.65, 369.46, 370.27, 371.09, 371.9, 372.72, 373.53, 374.34, 375.16, 375.97, 376.78, 377.6, 378.41, 379.22, 380.03, 380.84, 381.66, 382.47, 383.28, 384.09, 384.9, 385.71, 386.52, 387.34, 388.15, 388.96, 389.77, 390.58, 391.39, 392.2, 393.01, 393.82, 394.62, 395.43, 396.24, 397.05]"
$ws.Range("H19").Value2 = "[0.00770086, -0.00709508, 0.00182582, -0.00193501, -0.03777722, -0.00163547, -0.00249994, -0.00159035, -0.00860838, -0.005934, 0.00359084, 0.0053071, 0.00521964, 0.01838435, 0.0113439, 0.00018645, -0.00130239, -0.00392939, -0.02123891, -0.01861789, -0.00442397, 0.00221531, 0.00197348, 0.00433982, 0.00656391, 0.00772393, 0.00744643, 0.00721092, 0.00691338, 0.00660685, 0.00626031, 0.0059272, 0.00553537, 0.00515068, 0.00476844, 0.00441309, 0.00407636, 0.00372952, 0.00339784, 0.00308979, 0.00278495, 0.0025059, 0.00223471, 0.00198213, 0.00175449, 0.00153391, 0.00133415, 0.00115004, 0.00098961, 0.00085415, 0.00071921, 0.00061394, 0.00050784, 0.000422, 0.00035093, 0.00027316, 0.00021665, 0.00017891, 0.00014628, 0.00011272, 0.00008544, 0.00006722, 0.0000515, 0.00003649, 0.00002072, 0.00000832, 0.00000429, -0.0000021, -0.0000022, -0.00000356, -0.00000935, -0.00001333, -0.0000117, -0.00001569, -0.00002015, -0.00001611, -0.00000639, -0.00000348, -0.00001105, -0.00001074, -0.00000938, -0.00000603, -0.00000002, 0.0000088, 0.00001877, 0.00001588, 0.00001555, 0.00001249, 0.0000122, 0.00001424, 0.00000588, 0.00001866, 0.00002576, 0.00004831, 0.00003738, 0.00003228, 0.00005151, 0.00003754, 0.00004418, 0.00003845, 0.00004426, 0.00005321, 0.00005077, 0.00005776, 0.00006152, 0.00005503, 0.00006572, 0.00006447, 0.00007497, 0.00006871, 0.00007423, 0.00007492, 0.00008044, 0.00008543, 0.00010082, 0.00008981, 0.00008879, 0.00008732, 0.00008945, 0.00009184, 0.00009914, 0.00010118, 0.00011022, 0.00009938, 0.00010804, 0.00010921, 0.00011934, 0.00011225, 0.00012843, 0.00013692, 0.00014074, 0.00013363, 0.00014859, 0.00014807, 0.00014184, 0.00013851, 0.00014549, 0.00015069, 0.00014674, 0.00015669, 0.00015606, 0.00016257, 0.00017845, 0.00016934, 0.00016132, 0.00017124, 0.00018328, 0.00018204, 0.00018583, 0.00018105, 0.00017812, 0.00018305, 0.00018779, 0.00018992, 0.00019534, 0.00019682, 0.0002031, 0.00021514, 0.0002127, 0.00022214, 0.00020962, 0.00021476, 0.00021375, 0.0002249, 0.00022124, 0.00022573, 0.00023641, 0.00022327, 0.00023072, 0.00023193, 0.00023677, 0.00024845, 0.00025186, 0.00024623, 0.00025554, 0.0002411, 0.00025152, 0.00025919, 0.00025753, 0.0002492, 0.00026945, 0.00026774, 0.00027365, 0.00026696, 0.00027435, 0.00028988, 0.00028516, 0.00027722, 0.00028506, 0.00029847, 0.00031226, 0.00031051, 0.00030474, 0.00030534, 0.00030925, 0.00030493, 0.00030638, 0.00030454, 0.00029776, 0.00031279, 0.0002934, 0.00031017, 0.00032976, 0.00032801, 0.00035029, 0.00034698, 0.00033092, 0.00036255, 0.00032341, 0.00032763, 0.00034649, 0.00030852, 0.00034454, 0.00032961, 0.00034066, 0.00034117, 0.00034804, 0.00035471, 0.00037648, 0.00036621, 0.00038159, 0.00037974, 0.00036068, 0.00038973, 0.00037874, 0.00037441, 0.00038894, 0.00037116, 0.00037674, 0.00040048, 0.00039139, 0.00038325, 0.00041933, 0.00041024, 0.00043935, 0.00042598, 0.00043226, 0.00042031, 0.00042671, 0.00041124, 0.00038965, 0.00041458, 0.00043792, 0.00045666, 0.00042614, 0.00044786, 0.00044947, 0.00045128, 0.00043841, 0.00045139, 0.00043823, 0.00042743, 0.00037792, 0.00041029, 0.0004261, 0.00044434]"
$ws.Range("H20").Value2 = "[0.00148222, 0.00694155, 0.00425531, 0.00020583, 0.04311729, -0.00190723, 0.00008539, -0.00165836, 0.00151697, 0.00764164, -0.00154779, 0.00370061, 0.01886876, 0.03253636, 0.04243688, 0.05712787, 0.06995355, 0.08008836, 0.09710846, 0.1014758, 0.07692964, 0.07511929, 0.06789915, 0.06142195, 0.05256274, 0.0443899, 0.03736307, 0.03107416, 0.02550723, 0.02081273, 0.01678947, 0.01350332, 0.01078277, 0.0085821, 0.00680271, 0.00535675, 0.00421347, 0.00331048, 0.00261052, 0.00204349, 0.00161026, 0.00126565, 0.00099242, 0.00079095, 0.00062504, 0.00050219, 0.0004006, 0.0003355, 0.00028256, 0.00022531, 0.00019589, 0.00016879, 0.00014821, 0.00011631, 0.0000998, 0.00009875, 0.00007917, 0.00006016, 0.00004847, 0.00003939, 0.00003081, 0.00002348, 0.00000811, 0.00001479, 0.00001049, 0.00000755, 0.00000649, 0.00000303, -0.00000055, -0.00000116, -0.00000673, -0.00000468, -0.00001053, -0.00000481, -0.00000767, -0.00000687, -0.00000543, -0.00000781, -0.00000297, -0.00000735, -0.00000215, -0.00000628, -0.00001165, -0.00000287, -0.00000223, 0.00000728, 0.00000758, 0.00003454, 0.00002253, 0.00001535, 0.00002469, 0.00001791, 0.00002272, -0.00000035, 0.00000904, 0.00001782, 0.0000139, 0.00001438, 0.00002027, 0.00002845, 0.00002881, 0.00002804, 0.00003855, 0.00003822, 0.00003642, 0.0000473, 0.00004235, 0.00004054, 0.00003191, 0.00004811, 0.00004331, 0.0000545, 0.00005, 0.00004565, 0.00003864, 0.00005221, 0.0000561, 0.00006945, 0.00006677, 0.00005933, 0.00006624, 0.00007395, 0.00007711, 0.00008227, 0.00007404, 0.00008283, 0.00007569, 0.00008871, 0.00007202, 0.0000701, 0.00006785, 0.00008288, 0.00007065, 0.0000802, 0.00008602, 0.00009691, 0.00008831, 0.00009996, 0.00009736, 0.00009495, 0.00009612, 0.00010205, 0.00009163, 0.00010895, 0.00012551, 0.0001215, 0.00010939, 0.00012174, 0.00011344, 0.00012121, 0.00011066, 0.00011744, 0.00011828, 0.00012836, 0.0001261, 0.00013591, 0.00014283, 0.00014256, 0.00014432, 0.00014671, 0.00015578, 0.00013898, 0.00014799, 0.0001677, 0.00015912, 0.00017563, 0.00015877, 0.00016754, 0.0001748, 0.00017734, 0.0001847, 0.00019045, 0.00018323, 0.00019716, 0.00020503, 0.00021116, 0.00019593, 0.00021476, 0.00022371, 0.0002247, 0.00021386, 0.00024525, 0.00021305, 0.00023827, 0.00024244, 0.00025191, 0.00025414, 0.00028082, 0.00027033, 0.00027706, 0.00026555, 0.00025358, 0.00024676, 0.00026594, 0.00026636, 0.00026987, 0.00031013, 0.00030855, 0.00030851, 0.00032055, 0.00033466, 0.00034665, 0.00033285, 0.0003532, 0.00033182, 0.00033857, 0.00035741, 0.00031882, 0.00034114, 0.00035815, 0.00033801, 0.00036904, 0.00036752, 0.0003841, 0.0003714, 0.00037762, 0.00036031, 0.00039505, 0.00036625, 0.00041256, 0.00039857, 0.00039754, 0.00042607, 0.00040806, 0.00042922, 0.00044225, 0.0004252, 0.00041705, 0.00043993, 0.000421, 0.00043537, 0.00043257, 0.00043441, 0.00044707, 0.00041929, 0.00040509, 0.00043061, 0.00044794, 0.00045345, 0.00048101, 0.00049363, 0.00054348, 0.00051176, 0.00050844, 0.00048578, 0.00049046, 0.00050982, 0.00053094, 0.00050379, 0.00054719, 0.00051795, 0.00054938, 0.00054086, 0.00057717, 0.00062508, 0.0006358]"
$ws.Range("H21").Value2 = "[39, 34, 35, 30, 28, 31, 35, 26, 26, 32, 29, 46, 159, 412, 1219, 2915, 5582, 8814, 11845, 14148, 15672, 16509, 16861, 17016, 17172, 17368, 17783, 18402, 19303, 20468, 21923, 23764, 25908, 28399, 31138, 34116, 37132, 40065, 42722, 44828, 46221, 46844, 46585, 45659, 44239, 42474, 40614, 38939, 37452, 36237, 35385, 34824, 34569, 34548, 34840, 35376, 36082, 37019, 38167, 39510, 40947, 42575, 44290, 46079, 47778, 49277, 50696, 51716, 52243, 52227, 51708, 50703, 49175, 47332, 45321, 43038, 40851, 38741, 36744, 34964, 33342, 31908, 30673, 29639, 28761, 28039, 27483, 27123, 26859, 26743, 26795, 27008, 27336, 27809, 28397, 29111, 29942, 30811, 31762, 32724, 33638, 34400, 35122, 35627, 35880, 35919, 35696, 35205, 34533, 33678, 32667, 31647, 30526, 29414, 28364, 27379, 26461, 25642, 24912, 24279, 23781, 23372, 23072, 22852, 22723, 22704, 22742, 22876, 23059, 23305, 23603, 23986, 24375, 24803, 25339, 25881, 26424, 27034, 27641, 28264, 28881, 29482, 30065, 30617, 31153, 31611, 31999, 32348, 32573, 32716, 32704, 32648, 32475, 32179, 31787, 31325, 30790, 30137, 29477, 28775, 28048, 27288, 26535, 25846, 25075, 24358, 23604, 22899, 22186, 21473, 20784, 20161, 19555, 19021, 18503, 18032, 17603, 17239, 16904, 16582, 16291, 16037, 15738, 15499, 15252, 15007, 14776, 14565, 14337, 14129, 13931, 13741, 13586, 13474, 13370, 13257, 13197, 13109, 13053, 13015, 12980, 12950, 12929, 12895, 12894, 12874, 12876, 12844, 12817, 12803, 12771, 12719, 12693, 12640, 12575, 12507, 12452, 12413, 12309, 12240, 12048, 11844, 11635, 11407, 11201, 11019, 10808, 10622, 10465, 10290, 10117, 9961, 9845, 9731, 9545, 9282, 9050, 8838, 8638, 8440, 8239, 8093, 7918, 7795, 7673, 7546, 7400, 7287, 7211, 7198, 7116, 6903, 6496, 5991, 5467, 4895]"
$ws.Range("H22").Value2 = 217
$ws.Range("H23").Value2 = 240

# Apply the small Arial-10pt font used for the new calibration-name cells
# (column G, rows 17-23). Size/color already default to 10pt black here, so
# only the face needs to change.
for ($r = 17; $r -le 23; $r++) {
    $ws.Range("G$r").Font.Name = "Arial"
}

# New numeric coefficient format (2-decimal) for the lower/upper wavelength
# limits entered as plain numbers in H22:H23.
$ws.Range("H22:H23").NumberFormat = "0.00"

# I17:I23 should be blank (the old "No calibration coefficient" note from
# row 17 no longer applies now that real coefficients are present).
$ws.Range("I17:I23").ClearContents()

Write-Host "Edit complete"
